$d = $word.ActiveDocument

$d.Content.Find.Execute("BANK NAME - NationalBank", $true, $false, $false, $false, $false,
                         $true, 1, $false, "BANK NAME - BankUkraine", 2)

$d.Content.Find.Execute("CURRENCY NAME - RUB", $true, $false, $false, $false, $false,
                         $true, 1, $false, "CURRENCY NAME - USD", 2)

$d.Content.Find.Execute("DATE - 20.01.2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "DATE - 2019-09-02", 2)

$d.Content.Find.Execute("SALE RATE - 0.3948", $true, $false, $false, $false, $false,
                         $true, 1, $false, "SALE RATE - 999999.999", 2)
